$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "while city not in (...)" -> "while city not in CITY_DATA"
#    (paragraph that also still holds the <w:lastRenderedPageBreak/> marker)
# ---------------------------------------------------------------------------
$pWhile = $d.Paragraphs.Item(28)
$rWhile = $pWhile.Range.Duplicate
$foundWhile = $rWhile.Find.Execute( `
    "('chicago', 'new york city', 'washington')", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "CITY_DATA", 2)

# ---------------------------------------------------------------------------
# 2) "if city in (...)" -> "if city in CITY_DATA"
# ---------------------------------------------------------------------------
$pIf = $d.Paragraphs.Item(30)
$rIf = $pIf.Range.Duplicate
$foundIf = $rIf.Find.Execute( `
    "('chicago','new york city','washington')", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "CITY_DATA", 2)

# ---------------------------------------------------------------------------
# 3) Relocate the _GoBack bookmark: it used to sit right after "print('-'*60)"
#    further down the document; it should now sit inside "new york" (between
#    "yo" and "rk") on the "city = input(...)" line.
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$pInput = $d.Paragraphs.Item(29)
$rYork = $pInput.Range.Duplicate
$foundYork = $rYork.Find.Execute("york", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $rYork.Start + 2
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "while replace: $foundWhile"
Write-Output "if replace: $foundIf"
Write-Output "york found: $foundYork"
Write-Output $d.Paragraphs.Item(28).Range.Text
Write-Output $d.Paragraphs.Item(29).Range.Text
Write-Output $d.Paragraphs.Item(30).Range.Text
